$d = $word.ActiveDocument

# --- 1. Remove the stray "_GoBack" bookmark that sits after
#        "Good at studying, researching, collaborating, and documentation."
#        (it will be re-created at the new end-of-document edit point below).
$anchorText = "Good at studying, researching, collaborating, and documentation."
$findRange = $d.Content
$found = $findRange.Find.Execute($anchorText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $para = $findRange.Paragraphs(1)
    $pr = $para.Range
    $xml = $pr.WordOpenXML
    $idx = $xml.IndexOf($anchorText)
    if ($idx -ge 0) {
        $prefix = $xml.Substring(0, $idx)
        $pStart = $prefix.LastIndexOf("<w:p ")
        $pEndTag = $xml.IndexOf("</w:p>", $idx)
        $pEnd = $pEndTag + 6
        $paraXml = $xml.Substring($pStart, $pEnd - $pStart)
        $cleanXml = $paraXml -replace '<w:bookmarkStart[^/]*/>', '' -replace '<w:bookmarkEnd[^/]*/>', ''
        $pr.InsertXML($cleanXml)
    }
}

# --- 2. Append a new bullet describing the E-Sports club, in the same
#        "Accomplishments / Activities" list (numId 14), with the
#        "_GoBack" bookmark now marking this as the most-recent edit.
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:lastRenderedPageBreak/><w:t>Part of the E-Sports club for York County School of Technology. I play in Super Smash Bros. Ultimate and Mario Kart 8 Deluxe.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertXML($newParaXml)
